# 地图及关卡设计表格 - split the single default sheet into two design tables:
#   Sheet 1 "剧情线索"  (story-clue design table)
#   Sheet 2 "控制机关"  (control-mechanism design table)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename the existing sheet -----------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "剧情线索"

# --- Sheet 2: clone sheet 1 (keeps formatting/namespaces) then rename -----
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "控制机关"

# ===========================================================================
# 剧情线索 (sheet 1) contents
# ===========================================================================
$ws1.Range("A1").Value = "剧情线索编号"
$ws1.Range("B1").Value = "地图标注"
$ws1.Range("C1").Value = "物品名称"
$ws1.Range("D1").Value = "物品描述"
$ws1.Range("E1").Value = "获取方式"
$ws1.Range("F1").Value = "互动内容"

$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = "Clue 1#"
$ws1.Range("C2").Value = "破碎的笔记"
$ws1.Range("D2").Value = "这是一则笔记上撕下的一页，字迹已经模糊不堪了"
$ws1.Range("E2").Value = "调查房间的书柜"
$ws1.Range("F2").Value = "“这个世界已经到头，必须尽快找到那个地方...”"

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = "Clue 2#"
$ws1.Range("C3").Value = "房间一角的血迹"
$ws1.Range("D3").Value = "墙上的痕迹已经斑驳，显然已经过了很长时间，但是这是血迹没错"
$ws1.Range("E3").Value = "调查房间的墙壁"
$ws1.Range("F3").Value = "“墙上的痕迹已经斑驳，看来是很久以前留下的，但这是血迹没错...”"

$ws1.Range("A4").Value = 3
$ws1.Range("A5").Value = 4

# column widths (approximate Excel's own best-fit result as closely as
# the ColumnWidth setter's character-unit rounding allows)
$ws1.Columns.Item(1).ColumnWidth = 11.5
$ws1.Columns.Item(2).ColumnWidth = 11.5
$ws1.Columns.Item(3).ColumnWidth = 16.36
$ws1.Columns.Item(4).ColumnWidth = 60.07
$ws1.Columns.Item(5).ColumnWidth = 22.07
$ws1.Columns.Item(6).ColumnWidth = 65.5

# ===========================================================================
# 控制机关 (sheet 2) contents
# ===========================================================================
$ws2.Range("A1").Value = "控制机关编号"
$ws2.Range("B1").Value = "地图标注"
$ws2.Range("C1").Value = "机关形态"
$ws2.Range("D1").Value = "机关描述"
$ws2.Range("E1").Value = "需要物品"
$ws2.Range("F1").Value = "物品地图标注"
$ws2.Range("G1").Value = "破解条件"
$ws2.Range("H1").Value = "破解结果"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "Door 1#"
$ws2.Range("C2").Value = "紧锁的门"
$ws2.Range("D2").Value = "紧锁的门，找到钥匙才能打开"
$ws2.Range("E2").Value = "钥匙"
$ws2.Range("F2").Value = "Key 1#"
$ws2.Range("G2").Value = "在World 2#获取钥匙Key 1#"
$ws2.Range("H2").Value = "打开门"

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = "Door 2#"
$ws2.Range("C3").Value = "水底的门"
$ws2.Range("D3").Value = "角色无法进入水底，所以无法开门"
$ws2.Range("E3").Value = "NA"
$ws2.Range("F3").Value = "NA"
$ws2.Range("G3").Value = "进入World 3#"
$ws2.Range("H3").Value = "World 3#世界水已经干涸，可以直接跳下去开门"

$ws2.Range("A4").Value = 3
$ws2.Range("A5").Value = 4

$ws2.Columns.Item(1).ColumnWidth = 12.21
$ws2.Columns.Item(4).ColumnWidth = 30.93
$ws2.Columns.Item(5).ColumnWidth = 11.79
$ws2.Columns.Item(6).ColumnWidth = 13.79
$ws2.Columns.Item(7).ColumnWidth = 24.79
$ws2.Columns.Item(8).ColumnWidth = 43.93

# --- Selections -------------------------------------------------------------
# sheet1 -> C4 ; sheet2 -> H3 ; sheet2 ends up the active/visible tab since
# it is selected last (matches activeTab="1" / tabSelected="1" in the target)
$ws1.Range("C4").Select()
$ws2.Range("H3").Select()
